# Update the "question_answers" sheet (B2:B81) and the "outputs" sheet (B2:B9)
# with the new test-answer / score values.

$wb = $excel.ActiveWorkbook

$wsAnswers = $wb.Worksheets.Item("question_answers")
$wsOutputs = $wb.Worksheets.Item("outputs")

# New answers for question_answers!B<row> (row number => value). Only rows
# whose answer actually changed are listed here.
$answers = @{
    2  = "1"
    3  = "2"
    4  = "2"
    5  = "3"
    6  = "4"
    7  = "1"
    8  = "3"
    9  = "5"
    10 = "4"
    11 = "4"
    12 = "1"
    13 = "2"
    15 = "1"
    16 = "3"
    17 = "1"
    18 = "3"
    19 = "2"
    20 = "1"
    21 = "5"
    22 = "3"
    23 = "1"
    24 = "2"
    26 = "3"
    28 = "4"
    29 = "3"
    30 = "5"
    34 = "2"
    35 = "5"
    36 = "3"
    37 = "5"
    38 = "2"
    40 = "1"
    42 = "4"
    43 = "1"
    44 = "3"
    45 = "4"
    46 = "2"
    47 = "4"
    48 = "5"
    49 = "2"
    50 = "2"
    51 = "1"
    52 = "4"
    53 = "5"
    55 = "1"
    56 = "5"
    57 = "5"
    58 = "1"
    59 = "1"
    60 = "3"
    61 = "5"
    62 = "4"
    63 = "2"
    65 = "3"
    66 = "2"
    67 = "5"
    68 = "1"
    69 = "1"
    70 = "1"
    71 = "1"
    72 = "1"
    73 = "3"
    74 = "5"
    75 = "3"
    76 = "2"
    77 = "5"
    78 = "1"
    79 = "4"
    80 = "5"
    81 = "1"
}

foreach ($row in $answers.Keys) {
    # Force a text number format immediately before writing the value so
    # Excel stores the answer as text (matching the original inlineStr
    # cell type) instead of inferring a number.
    $cell = $wsAnswers.Range("B$row")
    $cell.NumberFormat = "@"
    $cell.Value = $answers[$row]
}

# New scores for outputs!B<row> (row number => value). Only rows whose
# score actually changed are listed here.
$outputs = @{
    2 = 27
    3 = 27
    4 = 24
    5 = 38
    6 = 29
    7 = 20
    9 = 30
}

foreach ($row in $outputs.Keys) {
    $wsOutputs.Range("B$row").Value = $outputs[$row]
}
